$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Unprotect()

# Update the confidential disclaimer date string (2021-05-11 -> 2021-05-12)
$ws.Range("A59").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + [char]10 + "Model holdings provided as of 2021-05-12 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for rows 2-56
$ws.Range("D2").Value = 0.01318122222338805
$ws.Range("E2").Value = -0.01042908224076278
$ws.Range("D3").Value = 0.01086121040450226
$ws.Range("E3").Value = -0.04946996466431097
$ws.Range("D4").Value = 0.01052522726695055
$ws.Range("E4").Value = -0.017911975435005
$ws.Range("D5").Value = 0.0113791002441152
$ws.Range("E5").Value = -0.01650871293182499
$ws.Range("D6").Value = 0.01099654629758421
$ws.Range("E6").Value = -0.004408523144746601
$ws.Range("D7").Value = 0.0141861417279891
$ws.Range("E7").Value = 0.009160305343511643
$ws.Range("D8").Value = 0.01113435100463615
$ws.Range("E8").Value = -0.03140495867768578
$ws.Range("D9").Value = 0.01123265468816831
$ws.Range("E9").Value = -0.03348785166240409
$ws.Range("D10").Value = 0.01055215978298676
$ws.Range("E10").Value = -0.01978050025523215
$ws.Range("D11").Value = 0.01101034921205277
$ws.Range("E11").Value = -0.002904754624675099
$ws.Range("D12").Value = 0.4377611156525442
$ws.Range("E12").Value = 0
$ws.Range("D13").Value = 0.01175155449713261
$ws.Range("E13").Value = -0.02694805194805194
$ws.Range("D14").Value = 0.01075550027906014
$ws.Range("E14").Value = -0.01416885772713994
$ws.Range("D15").Value = 0.01015670067252176
$ws.Range("E15").Value = -0.0262518230432669
$ws.Range("D16").Value = 0.009987250259127271
$ws.Range("E16").Value = -0.0242477359041775
$ws.Range("D17").Value = 0.009764159251294008
$ws.Range("E17").Value = -0.04108723135271808
$ws.Range("D18").Value = 0.008703242557100999
$ws.Range("E18").Value = -0.0462633451957295
$ws.Range("D19").Value = 0.009514809040325437
$ws.Range("E19").Value = -0.01738453554748309
$ws.Range("D20").Value = 0.01067649823202059
$ws.Range("E20").Value = -0.02102165230187092
$ws.Range("D21").Value = 0.01201661534245563
$ws.Range("E21").Value = -0.06639771390149596
$ws.Range("D22").Value = 0.01179520761687463
$ws.Range("E22").Value = 0.002540220152412953
$ws.Range("D23").Value = 0.01118372728403587
$ws.Range("E23").Value = -0.01791089704996984
$ws.Range("D24").Value = 0.01269194818206358
$ws.Range("E24").Value = -0.06427939876215738
$ws.Range("D25").Value = 0.01233868334672197
$ws.Range("E25").Value = -0.07275902211874274
$ws.Range("D26").Value = 0.01157155551495728
$ws.Range("E26").Value = -0.04387291981845698
$ws.Range("D27").Value = 0.01223925747502163
$ws.Range("E27").Value = -0.05314213412062418
$ws.Range("D28").Value = 0.01478325805227523
$ws.Range("E28").Value = -0.008767535070140386
$ws.Range("D29").Value = 0.01138235458980291
$ws.Range("E29").Value = -0.0002957704821060503
$ws.Range("D30").Value = 0.007133750185090908
$ws.Range("E30").Value = -0.001635991820040816
$ws.Range("D31").Value = 0.005080426384372003
$ws.Range("E31").Value = -0.0320834943950522
$ws.Range("D32").Value = 0.009538038335406668
$ws.Range("E32").Value = -0.02827225130890065
$ws.Range("D33").Value = 0.010810936374568
$ws.Range("E33").Value = 0.01023479831426855
$ws.Range("D34").Value = 0.01024670016360942
$ws.Range("E34").Value = -0.0005037783375315685
$ws.Range("D35").Value = 0.01007298543517587
$ws.Range("E35").Value = 0.01978565539983501
$ws.Range("D36").Value = 0.009768423566333072
$ws.Range("E36").Value = -0.005698005698005493
$ws.Range("D37").Value = 0.01086255703030407
$ws.Range("E37").Value = 0.005268703898840821
$ws.Range("D38").Value = 0.01143678071595942
$ws.Range("E38").Value = -0.0156502968159743
$ws.Range("D39").Value = 0.01450226213496411
$ws.Range("E39").Value = -0.03466633651108086
$ws.Range("D40").Value = 0.01085694608946319
$ws.Range("E40").Value = -0.03220738413197166
$ws.Range("D41").Value = 0.01291863019203501
$ws.Range("E41").Value = -0.02182070882557341
$ws.Range("D42").Value = 0.01152857570811616
$ws.Range("E42").Value = -0.03754392454225997
$ws.Range("D43").Value = 0.01146528429543107
$ws.Range("E43").Value = -0.02293259207783194
$ws.Range("D44").Value = 0.01051254654065017
$ws.Range("E44").Value = -0.007856616744414402
$ws.Range("D45").Value = 0.01155292719136557
$ws.Range("E45").Value = -0.02646915978630393
$ws.Range("D46").Value = 0.01112694456272619
$ws.Range("E46").Value = -0.04082538273796321
$ws.Range("D47").Value = 0.01023233615505678
$ws.Range("E47").Value = -0.02671579917088918
$ws.Range("D48").Value = 0.00946677938672753
$ws.Range("E48").Value = -0.02702702702702708
$ws.Range("D49").Value = 0.009764047032477189
$ws.Range("E49").Value = -0.05313243457573347
$ws.Range("D50").Value = 0.01024108922276855
$ws.Range("E50").Value = -0.06548323471400397
$ws.Range("D51").Value = 0.009284760465849484
$ws.Range("E51").Value = -0.01040634291377607
$ws.Range("D52").Value = 0.01031964239454082
$ws.Range("E52").Value = -0.04262722923009987
$ws.Range("D53").Value = 0.008607632125172456
$ws.Range("E53").Value = 0.002659574468085291
$ws.Range("D54").Value = 0.004351284622100047
$ws.Range("E54").Value = -0.03713733075435199
$ws.Range("D55").Value = 0.004185312992026908
$ws.Range("E55").Value = -0.01994851994851998
$ws.Range("D56").Value = 0.9999999999999997
$ws.Range("E56").Value = -0.01354767276970592

# Restore sheet protection (original password cannot be recovered from its hash,
# but re-applying protection keeps the sheet's protected state consistent)
$ws.Protect("")
